# Generate Report for Handoff
# Adds a new file's row ("e9181077-1948-439a-9c0a-51f40a39c601") to the
# Overview sheet and to each locale sheet (zh-cn, de-de), mirroring the
# existing b96bf4c9-... row, plus the matching hyperlinks.

$wb = $excel.ActiveWorkbook

$newGuid = "e9181077-1948-439a-9c0a-51f40a39c601"
$newHash = "7ccbed387454081580f8016f6bcdeaece521729f"

# ---------------------------------------------------------------------
# Overview sheet -> new row 3
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = ($newGuid + ".md")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-30-19 16:30:39"

$wsOverview.Range("A3").Style = "HyperLink"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e77b3cc790ca4da7a1c85b5aa04d80ed22b425ac/e2e/$newGuid.md",
    "",
    "",
    ($newGuid + ".md")
)

# ---------------------------------------------------------------------
# zh-cn sheet -> new row 3
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value = ($newGuid + ".md")
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-19 16:30:33"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value = "Include"

$wsZh.Range("A3").Style = "HyperLink"
$wsZh.Range("B3").Style = "HyperLink"
$wsZh.Range("D3").Style = "HyperLink"
$wsZh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e77b3cc790ca4da7a1c85b5aa04d80ed22b425ac/e2e/$newGuid.md",
    "",
    "",
    ($newGuid + ".md")
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e77b3cc790ca4da7a1c85b5aa04d80ed22b425ac/e2e/$newGuid.md",
    "",
    "",
    ".md"
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a9a9c941c371b1c86dc5bf70cf3026baba5061a4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newGuid.$newHash.zh-cn.xlf",
    "",
    "",
    "$newGuid.$newHash.zh-cn.xlf"
)

# ---------------------------------------------------------------------
# de-de sheet -> new row 3
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value = ($newGuid + ".md")
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-19 16:30:39"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value = "Include"

$wsDe.Range("A3").Style = "HyperLink"
$wsDe.Range("B3").Style = "HyperLink"
$wsDe.Range("D3").Style = "HyperLink"
$wsDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e77b3cc790ca4da7a1c85b5aa04d80ed22b425ac/e2e/$newGuid.md",
    "",
    "",
    ($newGuid + ".md")
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e77b3cc790ca4da7a1c85b5aa04d80ed22b425ac/e2e/$newGuid.md",
    "",
    "",
    ".md"
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3905c9c048c8079f1a89b20791eab1e86e3b77ba/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newGuid.$newHash.de-de.xlf",
    "",
    "",
    "$newGuid.$newHash.de-de.xlf"
)
